$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new "Register" worksheet right after "Login".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Register"

# Fill the cells in the order that matches the commit's shared-string layout
# (A1 reuses "john" from sheet1, the remaining new strings come from
#  row2 -> row3 -> B1 -> C column).
$ws2.Range("A1").Value = "john"
$ws2.Range("A2").Value = "rock"
$ws2.Range("B2").Value = "baby"
$ws2.Range("A3").Value = "justine"
$ws2.Range("B3").Value = "biber"
$ws2.Range("B1").Value = "rat"
$ws2.Range("C1").Value = 8870034785
$ws2.Range("C2").Value = 3223445666
$ws2.Range("C3").Value = 12112122

# Apply an explicit "General" number format across the new data so the
# style table grows a second cellXfs entry (applyNumberFormat).
$ws2.Range("A1:C3").NumberFormat = "General"

# Column C needs to be a bit wider to fit the phone numbers.
$ws2.Columns.Item(3).ColumnWidth = 10.1666666666667

# Selections: Login keeps a B1 selection and is no longer the active tab;
# Register becomes the active sheet with C4 selected.
$ws1.Range("B1").Select() | Out-Null
$ws2.Range("C4").Select() | Out-Null
